$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$dateFmt = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Range("G2").Value = "2016-09-02 21:15:42"
$wsOverview.Range("G2").NumberFormat = $dateFmt

$wsZhCn.Range("H2").Value = "2016-09-02 21:15:36"
$wsZhCn.Range("H2").NumberFormat = $dateFmt
$wsZhCn.Range("K2").Value = "2016-09-02 21:16:05"
$wsZhCn.Range("K2").NumberFormat = $dateFmt

$wsDeDe.Range("H2").Value = "2016-09-02 21:15:42"
$wsDeDe.Range("H2").NumberFormat = $dateFmt
$wsDeDe.Range("K2").Value = "2016-09-02 21:16:17"
$wsDeDe.Range("K2").NumberFormat = $dateFmt
